$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Organisation (column B) for the existing Carrefour product row (row 2),
# which previously had no party/org name.
$ws.Range("B2").Value = "Carrefour"

# Add a new product row (row 5) for a product imported with a GLN but no org/party name.
$ws.Range("A5").Value = 6270190128403
$ws.Range("C5").Value = "Some product without an org name"
$ws.Range("F5").Value = "Some product"

# Match the selection Excel leaves behind after this kind of edit.
$ws.Range("D5").Select()
